# Daily attendance processing - 2026-01-18 07:11:23
# Normalize the "Recorded By" (column G) author-order text so that the
# dnasr281@gmail.com / backup@backdoor.com addresses list before "System".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value2

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "System, system, backup@backdoor.com") {
        $cell.Value2 = "system, System, backup@backdoor.com"
    }
}
